$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

foreach ($addr in @("D5","D6","D8","D11","D13","D14","D18","D20","D21","D22","D25","D26","D28","D31","D33","D34","D37","D42","D44","D45","D51")) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '34.956.83'
$ws.Range("E2").Value = '  -0.53%  '
$ws.Range("D3").Value = '1.817.29'
$ws.Range("E3").Value = '  -0.48%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '230.25'
$ws.Range("E5").Value = '  -1.14%  '
$ws.Range("D6").Value = '0.615'
$ws.Range("E6").Value = '  +0.38%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("D8").Value = '40.11'
$ws.Range("E8").Value = '  -5.52%  '
$ws.Range("E9").Value = '  +4.00%  '
$ws.Range("E10").Value = '  -0.87%  '
$ws.Range("D11").Value = '0.0989'
$ws.Range("E11").Value = '  -1.54%  '
$ws.Range("D12").Value = '2.080.22'
$ws.Range("E12").Value = '  -0.59%  '
$ws.Range("D13").Value = '11.29'
$ws.Range("E13").Value = '  +1.11%  '
$ws.Range("D14").Value = '0.668'
$ws.Range("E14").Value = '  +0.40%  '
$ws.Range("D15").Value = '1.821.22'
$ws.Range("E15").Value = '  -0.28%  '
$ws.Range("E16").Value = '  -1.57%  '
$ws.Range("D17").Value = '34.949.12'
$ws.Range("E17").Value = '  -0.46%  '
$ws.Range("D18").Value = '69.61'
$ws.Range("E18").Value = '  -0.57%  '
$ws.Range("E19").Value = '  -1.01%  '
$ws.Range("D20").Value = '240.51'
$ws.Range("E20").Value = '  +0.17%  '
$ws.Range("D21").Value = '12.07'
$ws.Range("E21").Value = '  +2.01%  '
$ws.Range("D22").Value = '4.65'
$ws.Range("E22").Value = '  +1.07%  '
$ws.Range("E23").Value = '  +0.13%  '
$ws.Range("E24").Value = '  +1.87%  '
$ws.Range("D25").Value = '173.52'
$ws.Range("E25").Value = '  +1.01%  '
$ws.Range("D26").Value = '7.82'
$ws.Range("E26").Value = '  +0.58%  '
$ws.Range("E27").Value = '  +2.13%  '
$ws.Range("D28").Value = '17.32'
$ws.Range("E28").Value = '  -1.16%  '
$ws.Range("E29").Value = '  -5.72%  '
$ws.Range("E30").Value = '  -0.04%  '
$ws.Range("D31").Value = '3.99'
$ws.Range("E31").Value = '  +2.39%  '
$ws.Range("E32").Value = '  -1.22%  '
$ws.Range("D33").Value = '3.96'
$ws.Range("E33").Value = '  -1.10%  '
$ws.Range("D34").Value = '1.24'
$ws.Range("E34").Value = '  +11.46%  '
$ws.Range("E35").Value = '  +2.03%  '
$ws.Range("E36").Value = '  +1.67%  '
$ws.Range("D37").Value = '92.37'
$ws.Range("E37").Value = '  -1.22%  '
$ws.Range("E38").Value = '  +6.97%  '
$ws.Range("D39").Value = '1.339.73'
$ws.Range("E39").Value = '  +1.27%  '
$ws.Range("E40").Value = '  +0.02%  '
$ws.Range("E41").Value = '  -1.36%  '
$ws.Range("D42").Value = '14.49'
$ws.Range("E42").Value = '  -1.80%  '
$ws.Range("E43").Value = '  -3.28%  '
$ws.Range("D44").Value = '2.42'
$ws.Range("E44").Value = '  -1.18%  '
$ws.Range("D45").Value = '2.76'
$ws.Range("E45").Value = '  -1.00%  '
$ws.Range("E46").Value = '  +2.19%  '
$ws.Range("E47").Value = '  -0.51%  '
$ws.Range("D48").Value = '1.997.80'
$ws.Range("E48").Value = '  -0.32%  '
$ws.Range("E49").Value = '  +0.08%  '
$ws.Range("E50").Value = '  +3.20%  '
$ws.Range("D51").Value = '97.34'
$ws.Range("E51").Value = '  -3.15%  '
